$d = $word.ActiveDocument

# 1) "...Python modules such as pytest framework, python logging..."
#    -> capitalize "pytest" to "Pytest" (keep the rest of the sentence intact).
$d.Content.Find.Execute("pytest framework, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pytest framework, ", 2)

# 2) "...Redis, MinIO object store, Java." -> "...Redis, MinIO object store, Java, SQL."
#    Insert a new ", SQL" item into the skills list right after ", Java".
$d.Content.Find.Execute(", Java.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", Java, SQL.", 2)
